# Trade #26 closed at 2026-02-17 08:03:21 - unknown UNKNOWN +0.000%

$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B5").Value = -0.35
$summary.Range("B6").Value = 26
$summary.Range("B9").Value = 26.92

# --- Strategy Status sheet (MarketMaking row) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 26
$status.Range("G4").Value = 26.92

# --- New trade row data (Trade #26) ---
$tradeNum = 26
$tradeDate = "2026-02-17"
$tradeTime = "08:03:14"
$strategy = "MarketMaking"
$side = "DOWN"
$entryPrice = 0.74
$exitPrice = 0.74
$status2 = "CLOSED"
$pnlPct = 0
$pnlDollar = 0
$capitalAfter = 99.54000000000001
$entrySlippage = 0
$exitSlippage = 0
$confidence = 0.6
$entryReason = "Normal spread capture: 19600 bps"
$exitReason = "early_exit"
$duration = 0.13

function Add-TradeRow($ws) {
    $ws.Range("A27").Value = $tradeNum
    # Force the date-looking / time-looking strings to stay text (matching
    # how the existing rows store them as inline/shared strings, not as
    # Excel date/time serials).
    $ws.Range("B27").NumberFormat = "@"
    $ws.Range("B27").Value = $tradeDate
    $ws.Range("C27").NumberFormat = "@"
    $ws.Range("C27").Value = $tradeTime
    $ws.Range("D27").Value = $strategy
    $ws.Range("E27").Value = $side
    $ws.Range("F27").Value = $entryPrice
    $ws.Range("G27").Value = $exitPrice
    $ws.Range("H27").Value = $status2
    $ws.Range("I27").Value = $pnlPct
    $ws.Range("J27").Value = $pnlDollar
    $ws.Range("K27").Value = $capitalAfter
    $ws.Range("L27").Value = $entrySlippage
    $ws.Range("M27").Value = $exitSlippage
    $ws.Range("N27").Value = $confidence
    $ws.Range("O27").Value = $entryReason
    $ws.Range("P27").Value = $exitReason
    $ws.Range("Q27").Value = $duration

    # Drop the temporary text number-format so the new row's cells end up
    # with the same (default/general) style as the rest of the sheet.
    $ws.Range("A27:Q27").ClearFormats()
}

# --- All Trades sheet ---
$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $allTrades

# --- MarketMaking sheet ---
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $marketMaking
